$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the obsolete "x8" (outFunc=2) sub-blocks of 3 rows each.
# These are the 6 groups that disappear in the new experiment design
# (148*, 158*, 248*, 258*, 348*, 358*) - delete from the bottom up so
# earlier row numbers stay valid.
$ws.Rows("53:55").Delete()
$ws.Rows("44:46").Delete()
$ws.Rows("35:37").Delete()
$ws.Rows("26:28").Delete()
$ws.Rows("17:19").Delete()
$ws.Rows("8:10").Delete()

# --- Append the new "8" experiment group (isU2Y=0, isY2Y=1) as rows 38-49.
$newRows = @(
  @(38, 8469, 0, 1, 1, 0, 0),
  @(39, "846A", 0, 1, 1, 0, 1),
  @(40, "846B", 0, 1, 1, 0, 2),
  @(41, 8479, 0, 1, 1, 1, 0),
  @(42, "847A", 0, 1, 1, 1, 1),
  @(43, "847B", 0, 1, 1, 1, 2),
  @(44, 8569, 0, 1, 2, 0, 0),
  @(45, "856A", 0, 1, 2, 0, 1),
  @(46, "856B", 0, 1, 2, 0, 2),
  @(47, 8579, 0, 1, 2, 1, 0),
  @(48, "857A", 0, 1, 2, 1, 1),
  @(49, "857B", 0, 1, 2, 1, 2)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
}

# Column A in this sheet is styled right-aligned; match that for the new
# labels except row 38, whose A cell was typed in directly without the
# style carrying over (matches source workbook).
$ws.Range("A39:A49").HorizontalAlignment = -4152
